function Copy-RowData {
    # Copies one full A:J row from $srcSheet/$srcRow to $dstSheet/$dstRow.
    # Style comes from Range.Copy (reliable); values are re-applied explicitly
    # afterwards because Copy() silently ignores empty-string source cells
    # and would otherwise leave the destination's previous value behind.
    param($srcSheet, $srcRow, $dstSheet, $dstRow)

    $srcRange = $srcSheet.Range("A" + $srcRow + ":J" + $srcRow)
    $dstRange = $dstSheet.Range("A" + $dstRow + ":J" + $dstRow)
    $srcRange.Copy($dstRange)

    $cols = @("A","B","C","D","E","F","G","H","I","J")
    foreach ($c in $cols) {
        $v = $srcSheet.Range($c + $srcRow).Value()
        if ($v -is [string]) {
            # Force literal text (column already carries quotePrefix styling
            # where needed) so numeric-looking strings like "1.250" aren't
            # silently reinterpreted as the number 1.25.
            $dstSheet.Range($c + $dstRow).Value = "'" + $v
        } else {
            $dstSheet.Range($c + $dstRow).Value = $v
        }
    }
}

$wb = $excel.ActiveWorkbook

$panels = $wb.Worksheets.Item("Add Panels")
$sheet1 = $wb.Worksheets.Item("Sheet1")

# --- Relocate the rows that now belong on the "more panel test data" sheet ---
# (Row numbers below refer to the original "Add Panels" layout, read before
# anything is overwritten.)
Copy-RowData $panels 11 $sheet1 8
Copy-RowData $panels 12 $sheet1 9
Copy-RowData $panels 13 $sheet1 10
Copy-RowData $panels 14 $sheet1 11
Copy-RowData $panels 9  $sheet1 12

# --- Pull row 10 up into row 9 on "Add Panels" ---
# (row 9's original content has already been relocated above, so it is now
# safe to overwrite.)
Copy-RowData $panels 10 $panels 9

# --- Remove the now-duplicated trailing rows from "Add Panels" ---
$panels.Rows("10:14").Delete()

# --- Rename the second worksheet ---
$sheet1.Name = "more panel test data"

# --- Update the loading-details column header text everywhere ---
$panels.Cells.Replace("40V (A)", "40V Rail(A)")
$sheet1.Cells.Replace("40V (A)", "40V Rail(A)")

# --- Restore each sheet's selection, leaving "Add Panels" as the active tab ---
$sheet1.Range("C29").Select()
$panels.Activate()
$panels.Range("A9").Select()
